# Linear calibration sheet was recorded with the roles of X and Y reversed:
# the program actually receives the (known) Y value and the (measured)
# sensor/X value, so the transfer function has to be inverted. Rename the
# first header from "xValue" to "yValue" - "sensorValue" in B1 is unaffected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "yValue"

# Leftover UI state from the edit session (selection moved to B10:B11).
$ws.Range("B10:B11").Select() | Out-Null
